# Updates cryptos price/volume data per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.059.61"
$ws.Range("E2").Value = "  +3.87%  "

# Row 3
$ws.Range("D3").Value = "1.683.70"
$ws.Range("E3").Value = "  +3.20%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'220.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.39%  "

# Row 6
$ws.Range("E6").Value = "  +2.02%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "'29.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.71%  "

# Row 9
$ws.Range("E9").Value = "  +2.64%  "

# Row 10
$ws.Range("D10").Value = "'0.0640"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.09%  "

# Row 11
$ws.Range("E11").Value = "  +0.76%  "

# Row 12
$ws.Range("D12").Value = "1.928.58"
$ws.Range("E12").Value = "  +3.39%  "

# Row 13
$ws.Range("D13").Value = "1.683.98"
$ws.Range("E13").Value = "  +3.23%  "

# Row 14
$ws.Range("D14").Value = "'10.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.92%  "

# Row 15
$ws.Range("D15").Value = "'0.607"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.86%  "

# Row 16
$ws.Range("D16").Value = "'4.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.99%  "

# Row 17
$ws.Range("D17").Value = "31.079.12"
$ws.Range("E17").Value = "  +3.93%  "

# Row 18
$ws.Range("D18").Value = "'66.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.04%  "

# Row 19
$ws.Range("D19").Value = "'247.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.79%  "

# Row 20
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
$ws.Range("D21").Value = "'0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("E22").Value = "  +3.72%  "

# Row 23
$ws.Range("D23").Value = "'10.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.90%  "

# Row 24
$ws.Range("E24").Value = "  -0.69%  "

# Row 25
$ws.Range("D25").Value = "'158.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.57%  "

# Row 26
$ws.Range("E26").Value = "  +2.80%  "

# Row 27
$ws.Range("E27").Value = "  +2.67%  "

# Row 28
$ws.Range("E28").Value = "  +1.33%  "

# Row 29
$ws.Range("E29").Value = "  -0.04%  "

# Row 30
$ws.Range("D30").Value = "'0.0499"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.18%  "

# Row 31
$ws.Range("E31").Value = "  +4.02%  "

# Row 32
$ws.Range("E32").Value = "  +3.57%  "

# Row 33
$ws.Range("E33").Value = "  +5.18%  "

# Row 34
$ws.Range("D34").Value = "1.514.62"
$ws.Range("E34").Value = "  +6.44%  "

# Row 35
$ws.Range("E35").Value = "  +2.96%  "

# Row 36
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.91%  "

# Row 37
$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "'83.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.48%  "

# Row 38
$ws.Range("D38").Value = "'0.611"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.70%  "

# Row 39
$ws.Range("E39").Value = "  +5.31%  "

# Row 40
$ws.Range("D40").Value = "'2.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.81%  "

# Row 41
$ws.Range("E41").Value = "  +0.25%  "

# Row 42
$ws.Range("E42").Value = "  +2.10%  "

# Row 43
$ws.Range("D43").Value = "'0.840"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

# Row 44
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("E45").Value = "  +2.38%  "

# Row 46
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("D47").Value = "'52.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.37%  "

# Row 48
$ws.Range("D48").Value = "'5.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.39%  "

# Row 49
$ws.Range("D49").Value = "1.819.00"
$ws.Range("E49").Value = "  +2.62%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0120"
$ws.Range("E50").Value = "  +8.86%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'93.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
